# Auto-generated edit script: updates Leve profit-tracking values
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# per the scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 425
$ws.Range("I6").Value = 300
$ws.Range("K6").Value = 900
$ws.Range("M6").Value = -788
$ws.Range("H9").Value = 15111177
$ws.Range("I9").Value = 13888972
$ws.Range("J9").Value = 20000000
$ws.Range("K9").Value = 13888972
$ws.Range("L9").Value = 20000000
$ws.Range("M9").Value = -13888803
$ws.Range("N9").Value = -20000338
$ws.Range("H64").Value = 2832.4675
$ws.Range("I64").Value = 2750
$ws.Range("J64").Value = 2921.6216
$ws.Range("K64").Value = 2750
$ws.Range("L64").Value = 2921.6216
$ws.Range("M64").Value = -2502
$ws.Range("N64").Value = -3417.6216
$ws.Range("H67").Value = 2832.4675
$ws.Range("I67").Value = 2750
$ws.Range("J67").Value = 2921.6216
$ws.Range("K67").Value = 2750
$ws.Range("L67").Value = 2921.6216
$ws.Range("M67").Value = -1892
$ws.Range("N67").Value = -4637.6216
$ws.Range("H106").Value = 1131.6666
$ws.Range("I106").Value = 1131.6666
$ws.Range("K106").Value = 1131.6666
$ws.Range("M106").Value = -500.6666
$ws.Range("H113").Value = 4548.846
$ws.Range("I113").Value = 2400
$ws.Range("J113").Value = 5891.875
$ws.Range("K113").Value = 2400
$ws.Range("L113").Value = 5891.875
$ws.Range("M113").Value = 854
$ws.Range("N113").Value = -12399.875
$ws.Range("H116").Value = 5551.1333
$ws.Range("J116").Value = 5298.875
$ws.Range("L116").Value = 5298.875
$ws.Range("N116").Value = -12182.875
$ws.Range("H123").Value = 48078.75
$ws.Range("J123").Value = 48078.75
$ws.Range("L123").Value = 48078.75
$ws.Range("N123").Value = -57878.75
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H130").Value = 37660
$ws.Range("J130").Value = 37660
$ws.Range("L130").Value = 37660
$ws.Range("N130").Value = -47700
$ws.Range("H134").Value = 44302
$ws.Range("J134").Value = 44302
$ws.Range("L134").Value = 44302
$ws.Range("N134").Value = -54442

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16236.55
$ws.Range("I32").Value = 4490.9165
$ws.Range("J32").Value = 33855
$ws.Range("K32").Value = 4490.9165
$ws.Range("L32").Value = 33855
$ws.Range("M32").Value = -4203.9165
$ws.Range("N32").Value = -34429
$ws.Range("H45").Value = 3227.84
$ws.Range("I45").Value = 2490.3333
$ws.Range("J45").Value = 5124.2856
$ws.Range("K45").Value = 2490.3333
$ws.Range("L45").Value = 5124.2856
$ws.Range("M45").Value = -2113.3333
$ws.Range("N45").Value = -5878.2856
$ws.Range("H123").Value = 2500000
$ws.Range("J123").Value = 2500000
$ws.Range("L123").Value = 2500000
$ws.Range("N123").Value = -2509800

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3565.45
$ws.Range("J105").Value = 7800
$ws.Range("L105").Value = 7800
$ws.Range("N105").Value = -11294
$ws.Range("H122").Value = 33992.5
$ws.Range("J122").Value = 33992.5
$ws.Range("L122").Value = 33992.5
$ws.Range("N122").Value = -43792.5
$ws.Range("H126").Value = 23254.285
$ws.Range("J126").Value = 23254.285
$ws.Range("L126").Value = 23254.285
$ws.Range("N126").Value = -33134.285
$ws.Range("H130").Value = 40730
$ws.Range("J130").Value = 40730
$ws.Range("L130").Value = 40730
$ws.Range("N130").Value = -50770
$ws.Range("H132").Value = 48875
$ws.Range("J132").Value = 48875
$ws.Range("L132").Value = 48875
$ws.Range("N132").Value = -58995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49249.5
$ws.Range("J20").Value = 49249.5
$ws.Range("L20").Value = 49249.5
$ws.Range("N20").Value = -49721.5
$ws.Range("H30").Value = 49249.5
$ws.Range("J30").Value = 49249.5
$ws.Range("L30").Value = 49249.5
$ws.Range("N30").Value = -49431.5
$ws.Range("H53").Value = 25675
$ws.Range("J53").Value = 25675
$ws.Range("L53").Value = 25675
$ws.Range("N53").Value = -26889
$ws.Range("H62").Value = 3827.2727
$ws.Range("I62").Value = 3746.6667
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 3746.6667
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -3122.6667
$ws.Range("N62").Value = -5248
$ws.Range("H65").Value = 3827.2727
$ws.Range("I65").Value = 3746.6667
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 18733.3335
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -15613.3335
$ws.Range("N65").Value = -26240
$ws.Range("H105").Value = 2134.0588
$ws.Range("I105").Value = 1057.0714
$ws.Range("J105").Value = 7160
$ws.Range("K105").Value = 1057.0714
$ws.Range("L105").Value = 7160
$ws.Range("M105").Value = 689.9286
$ws.Range("N105").Value = -10654
$ws.Range("H107").Value = 872.96
$ws.Range("I107").Value = 1215.0714
$ws.Range("J107").Value = 437.54544
$ws.Range("K107").Value = 1215.0714
$ws.Range("L107").Value = 437.54544
$ws.Range("M107").Value = 704.9286
$ws.Range("N107").Value = -4277.54544
$ws.Range("H108").Value = 15144.667
$ws.Range("I108").Value = 8290
$ws.Range("J108").Value = 18572
$ws.Range("K108").Value = 8290
$ws.Range("L108").Value = 18572
$ws.Range("M108").Value = -4450
$ws.Range("N108").Value = -26252
$ws.Range("H109").Value = 10942.857
$ws.Range("J109").Value = 10942.857
$ws.Range("L109").Value = 10942.857
$ws.Range("N109").Value = -13022.857
$ws.Range("H122").Value = 2700
$ws.Range("I122").Value = 1575
$ws.Range("K122").Value = 4725
$ws.Range("M122").Value = -2275
$ws.Range("H124").Value = 44375
$ws.Range("J124").Value = 44375
$ws.Range("L124").Value = 44375
$ws.Range("N124").Value = -49285
$ws.Range("H127").Value = 54960
$ws.Range("J127").Value = 54940
$ws.Range("L127").Value = 54940
$ws.Range("N127").Value = -64860
$ws.Range("H128").Value = 49249.5
$ws.Range("J128").Value = 49249.5
$ws.Range("L128").Value = 49249.5
$ws.Range("N128").Value = -59209.5
$ws.Range("H135").Value = 53512.5
$ws.Range("J135").Value = 53512.5
$ws.Range("L135").Value = 53512.5
$ws.Range("N135").Value = -63652.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 10431.111
$ws.Range("J46").Value = 15576
$ws.Range("L46").Value = 15576
$ws.Range("N46").Value = -15888
$ws.Range("H102").Value = 2574
$ws.Range("I102").Value = 1432
$ws.Range("J102").Value = 6000
$ws.Range("K102").Value = 1432
$ws.Range("L102").Value = 6000
$ws.Range("M102").Value = 190
$ws.Range("N102").Value = -9244
$ws.Range("H119").Value = 40252.332
$ws.Range("J119").Value = 40252.332
$ws.Range("L119").Value = 40252.332
$ws.Range("N119").Value = -49928.332
$ws.Range("H124").Value = 56500
$ws.Range("J124").Value = 56500
$ws.Range("L124").Value = 56500
$ws.Range("N124").Value = -66320
$ws.Range("H126").Value = 13218.32
$ws.Range("I126").Value = 3388.9092
$ws.Range("K126").Value = 10166.7276
$ws.Range("M126").Value = -7696.7276
$ws.Range("H133").Value = 20186.666
$ws.Range("J133").Value = 20186.666
$ws.Range("L133").Value = 20186.666
$ws.Range("N133").Value = -30306.666
$ws.Range("H135").Value = 47473.332
$ws.Range("J135").Value = 47473.332
$ws.Range("L135").Value = 47473.332
$ws.Range("N135").Value = -57613.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H54").Value = 38000
$ws.Range("J54").Value = 38000
$ws.Range("L54").Value = 38000
$ws.Range("N54").Value = -39288
$ws.Range("H74").Value = 166693000
$ws.Range("J74").Value = 166693000
$ws.Range("L74").Value = 166693000
$ws.Range("N74").Value = -166694996
$ws.Range("H77").Value = 166693000
$ws.Range("J77").Value = 166693000
$ws.Range("L77").Value = 500079000
$ws.Range("N77").Value = -500088984
$ws.Range("H80").Value = 35128
$ws.Range("J80").Value = 35128
$ws.Range("L80").Value = 35128
$ws.Range("N80").Value = -37374
$ws.Range("H83").Value = 35128
$ws.Range("J83").Value = 35128
$ws.Range("L83").Value = 105384
$ws.Range("N83").Value = -116616
$ws.Range("H92").Value = 39289
$ws.Range("J92").Value = 39289
$ws.Range("L92").Value = 39289
$ws.Range("N92").Value = -44281
$ws.Range("H123").Value = 40422.25
$ws.Range("J123").Value = 40422.25
$ws.Range("L123").Value = 40422.25
$ws.Range("N123").Value = -50222.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 29993.334
$ws.Range("J93").Value = 29993.334
$ws.Range("L93").Value = 29993.334
$ws.Range("N93").Value = -34985.334
$ws.Range("H109").Value = 18055.666
$ws.Range("J109").Value = 18055.666
$ws.Range("L109").Value = 18055.666
$ws.Range("N109").Value = -20829.666

